$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "Átlagos Kereset"

# Update average gross earnings values (column B) for the corresponding years
$ws.Range("B2").Value = 202525
$ws.Range("B3").Value = 213094
$ws.Range("B4").Value = 223060
$ws.Range("B5").Value = 230714
$ws.Range("B6").Value = 237695
$ws.Range("B7").Value = 247924
$ws.Range("B8").Value = 263171
$ws.Range("B9").Value = 297017
$ws.Range("B10").Value = 329943
$ws.Range("B11").Value = 367833
$ws.Range("B14").Value = 515766
$ws.Range("B15").Value = 589114
